$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Plain text/value updates (safe from numeric auto-conversion) ---
$ws.Range("D2").Value = '43.818.69'
$ws.Range("E2").Value = '  +0.94%  '
$ws.Range("D3").Value = '2.351.21'
$ws.Range("E3").Value = '  +0.25%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("E5").Value = '  +3.44%  '
$ws.Range("E6").Value = '  +1.92%  '
$ws.Range("E7").Value = '  +10.55%  '
$ws.Range("E8").Value = '  -0.08%  '
$ws.Range("E9").Value = '  +19.05%  '
$ws.Range("E10").Value = '  +2.96%  '
$ws.Range("E11").Value = '  +6.21%  '
$ws.Range("E12").Value = '  +2.43%  '
$ws.Range("D13").Value = '2.700.08'
$ws.Range("E13").Value = '  +0.18%  '
$ws.Range("E14").Value = '  +8.10%  '
$ws.Range("E15").Value = '  +6.11%  '
$ws.Range("E16").Value = '  +5.29%  '
$ws.Range("D17").Value = '2.359.48'
$ws.Range("E17").Value = '  +0.67%  '
$ws.Range("D18").Value = '43.848.76'
$ws.Range("E18").Value = '  +1.13%  '
$ws.Range("E19").Value = '  +2.86%  '
$ws.Range("E20").Value = '  +4.97%  '
$ws.Range("E21").Value = '  +3.29%  '
$ws.Range("E22").Value = '  +1.90%  '
$ws.Range("E23").Value = '  -0.01%  '
$ws.Range("E24").Value = '  -2.58%  '
$ws.Range("E25").Value = '  +2.96%  '
$ws.Range("E26").Value = '  +6.31%  '
$ws.Range("E27").Value = '  -0.10%  '
$ws.Range("E28").Value = '  +0.90%  '
$ws.Range("E29").Value = '  -1.24%  '
$ws.Range("E30").Value = '  +6.16%  '
$ws.Range("E31").Value = '  +1.74%  '
$ws.Range("E32").Value = '  +4.72%  '
$ws.Range("E33").Value = '  +3.69%  '
$ws.Range("E34").Value = '  +3.55%  '
$ws.Range("E35").Value = '  +4.19%  '
$ws.Range("E36").Value = '  +11.89%  '
$ws.Range("E37").Value = '  -4.27%  '
$ws.Range("E38").Value = '  -1.00%  '
$ws.Range("E39").Value = '  +6.21%  '
$ws.Range("E40").Value = '  +9.12%  '
$ws.Range("E41").Value = '  -0.06%  '
$ws.Range("E42").Value = '  -2.21%  '
$ws.Range("E43").Value = '  +2.90%  '
$ws.Range("E44").Value = '  +3.54%  '
$ws.Range("E45").Value = '  -0.64%  '
$ws.Range("E46").Value = '  +1.54%  '
$ws.Range("B47").Value = 'Aave'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("E47").Value = '  -1.45%  '
$ws.Range("B48").Value = 'Algorand'
$ws.Range("C48").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("E48").Value = '  +11.19%  '
$ws.Range("B49").Value = 'NEARProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("E49").Value = '  +2.37%  '
$ws.Range("B50").Value = 'Maker'
$ws.Range("C50").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D50").Value = '1.434.08'
$ws.Range("E50").Value = '  -0.20%  '
$ws.Range("E51").Value = '  +1.50%  '

# --- Numeric-looking price values: force text entry, then restore default style ---
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '237.00'
$ws.Range("D6").NumberFormat = "General"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '72.78'
$ws.Range("D7").NumberFormat = "General"
$ws.Range("D7").Style = "Normal"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.539'
$ws.Range("D9").NumberFormat = "General"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0988'
$ws.Range("D10").NumberFormat = "General"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '28.45'
$ws.Range("D11").NumberFormat = "General"
$ws.Range("D11").Style = "Normal"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '16.65'
$ws.Range("D14").NumberFormat = "General"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.64'
$ws.Range("D15").NumberFormat = "General"
$ws.Range("D15").Style = "Normal"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '77.85'
$ws.Range("D20").NumberFormat = "General"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.43'
$ws.Range("D21").NumberFormat = "General"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '254.00'
$ws.Range("D22").NumberFormat = "General"
$ws.Range("D22").Style = "Normal"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.51'
$ws.Range("D25").NumberFormat = "General"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '10.54'
$ws.Range("D26").NumberFormat = "General"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.28'
$ws.Range("D27").NumberFormat = "General"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '22.37'
$ws.Range("D28").NumberFormat = "General"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '172.68'
$ws.Range("D29").NumberFormat = "General"
$ws.Range("D29").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.131'
$ws.Range("D32").NumberFormat = "General"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.16'
$ws.Range("D33").NumberFormat = "General"
$ws.Range("D33").Style = "Normal"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.16'
$ws.Range("D35").NumberFormat = "General"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.04'
$ws.Range("D36").NumberFormat = "General"
$ws.Range("D36").Style = "Normal"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.40'
$ws.Range("D38").NumberFormat = "General"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0267'
$ws.Range("D39").NumberFormat = "General"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '19.70'
$ws.Range("D40").NumberFormat = "General"
$ws.Range("D40").Style = "Normal"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '8.81'
$ws.Range("D42").NumberFormat = "General"
$ws.Range("D42").Style = "Normal"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0980'
$ws.Range("D44").NumberFormat = "General"
$ws.Range("D44").Style = "Normal"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '4.44'
$ws.Range("D46").NumberFormat = "General"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '97.64'
$ws.Range("D47").NumberFormat = "General"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.181'
$ws.Range("D48").NumberFormat = "General"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.31'
$ws.Range("D49").NumberFormat = "General"
$ws.Range("D49").Style = "Normal"
